# New crime data collected - weekly CompStat update
# Updates report volume/number, week-covering dates, and the crime-complaint
# statistics table (rows 16-28) to reflect newly collected data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (shared-string rich-text runs)
# ---------------------------------------------------------------------------
# "Volume 31   Number  21" -> "...Number  22"
$ws.Range("A8").Replace("21", "22")

# "Report Covering the Week  5/20/2024  Through  5/26/2024"
#   -> "...5/27/2024  Through  6/2/2024"
$ws.Range("C9").Replace("5/20/2024", "5/27/2024")
$ws.Range("C9").Replace("5/26/2024", "6/2/2024")

# ---------------------------------------------------------------------------
# Crime statistics table updates (rows 16-28)
# ---------------------------------------------------------------------------

# --- Row 16 (Robbery) ---
# D16 and E16 switch from numeric (1 / -100) to the "not applicable" text
# placeholders ("0" / "***.*"), matching the style used in row 20.
$ws.Range("G20").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("G20").Copy()
$ws.Range("D16").PasteSpecial(-4163)

$ws.Range("H20").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("H20").Copy()
$ws.Range("E16").PasteSpecial(-4163)

$ws.Range("F16").Value = 1
$ws.Range("H16").Value = 0
$ws.Range("M16").Value = 122.222222222222
$ws.Range("N16").Value = -64.912280701754

# --- Row 17 (Fel. Assault) ---
# C17 switches from numeric (1) to the "0" text placeholder.
$ws.Range("G20").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("G20").Copy()
$ws.Range("C17").PasteSpecial(-4163)

$ws.Range("N17").Value = -68.421052631578

# --- Row 18 (Burglary) ---
$ws.Range("G18").Value = 2
$ws.Range("J18").Value = 3

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 11
$ws.Range("G19").Value = 6
$ws.Range("H19").Value = 83.333333333333
$ws.Range("I19").Value = 23
$ws.Range("J19").Value = 14
$ws.Range("K19").Value = 64.285714285714
$ws.Range("L19").Value = 228.571428571429
$ws.Range("M19").Value = -4.166666666666
$ws.Range("N19").Value = -57.407407407407

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 14
$ws.Range("G21").Value = 9
$ws.Range("H21").Value = 55.555555555555
$ws.Range("I21").Value = 50
$ws.Range("J21").Value = 25
$ws.Range("K21").Value = 100
$ws.Range("L21").Value = 138.095238095238
$ws.Range("M21").Value = 38.888888888888
$ws.Range("N21").Value = -66.442953020134

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 2

# D24 and E24 switch from numeric (1 / 200) to the text placeholders.
$ws.Range("G20").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("G20").Copy()
$ws.Range("D24").PasteSpecial(-4163)

$ws.Range("H20").Copy()
$ws.Range("E24").PasteSpecial(-4122)
$ws.Range("H20").Copy()
$ws.Range("E24").PasteSpecial(-4163)

$ws.Range("F24").Value = 7
$ws.Range("H24").Value = 250
$ws.Range("I24").Value = 12
$ws.Range("K24").Value = -7.692307692307
$ws.Range("L24").Value = 33.333333333333
$ws.Range("M24").Value = -50

# --- Row 25 (Retail Theft) ---
# C25 switches from the "0" text placeholder to a plain numeric value (1).
$ws.Range("F25").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C25").Value = 1

$ws.Range("F25").Value = 2
$ws.Range("I25").Value = 3

# --- Row 26 (Misd. Assault) ---
# C26 switches from numeric (1) to the "0" text placeholder.
$ws.Range("G20").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("G20").Copy()
$ws.Range("C26").PasteSpecial(-4163)

$ws.Range("D26").Value = 1
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 14
$ws.Range("K26").Value = -28.571428571428
$ws.Range("L26").Value = -16.666666666666

# --- Row 28 (Shooting Inc.) ---
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 0
